# Add a new second paragraph:
#   "Adding 2" + "nd" (superscript) + " line for test."
$d = $word.ActiveDocument

# Append a new paragraph right after the existing (only) paragraph.
$firstParaRange = $d.Paragraphs(1).Range
$firstParaRange.InsertParagraphAfter()

# Grab the freshly created (now second) paragraph and fill it with the
# full sentence first - splitting/formatting runs afterwards.
$secondPara = $d.Paragraphs(2)
$secondRange = $secondPara.Range
$secondRange.Collapse(0)
$secondRange.InsertAfter("Adding 2nd line for test.")

# Locate the "nd" substring within that paragraph and mark it superscript,
# matching Word's normal "ordinal suffix" autoformat run-splitting.
$ordinalRange = $secondPara.Range
$ordinalRange.Find.Execute("nd", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
$ordinalRange.Font.Superscript = $true
